$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped from
# 45202 (2023-10-03) to 45203 (2023-10-04) for every data row (rows 2-367).
$ws.Range("C2:C367").Value = 45203
